$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 105
$ws.Range("A105").Value = "post-2.jpg"
$ws.Range("B105").Value = 1
$ws.Range("C105").Value = ""
$ws.Range("D105").Value = "Scheduled Image only at None"
$ws.Range("E105").Value = "2025-07-28 10:25:05"
$ws.Range("F105").Value = "Uncategorized"

# Row 106
$ws.Range("A106").Value = "post-1.jpg"
$ws.Range("B106").Value = 2
$ws.Range("C106").Value = ""
$ws.Range("D106").Value = "Scheduled Image only at None"
$ws.Range("E106").Value = "2025-07-28 10:30:04"
$ws.Range("F106").Value = "Uncategorized"

# Row 107
$ws.Range("A107").Value = "post-2.jpg"
$ws.Range("B107").Value = 1
$ws.Range("C107").Value = "🔥 **BLOCKBUSTER DEALS OF THE DAY** 🔥`n💥 Price Crash Store`n⚡️ Up to 5% off`n👉🏻 amzaff.in/l3swo0g`n🌟 Ki"
$ws.Range("D107").Value = "Scheduled Image + Text at None"
$ws.Range("E107").Value = "2025-07-28 10:32:02"
$ws.Range("F107").Value = "Kid's Carnival"

# Row 108
$ws.Range("A108").Value = "post-1.jpg"
$ws.Range("B108").Value = 2
$ws.Range("C108").Value = "💧 **Everyday Essentials for Skin, Hair & Fragrance** 💧`n💥 Maximise earnings with Beauty commissions "
$ws.Range("D108").Value = "Scheduled Image + Text at None"
$ws.Range("E108").Value = "2025-07-28 10:35:03"
$ws.Range("F108").Value = "Daily Essentials"

# Row 109
$ws.Range("A109").Value = "post-1.jpg"
$ws.Range("B109").Value = 1
$ws.Range("C109").Value = "🔥 **BLOCKBUSTER DEALS OF THE DAY** 🔥`n💥 Price Crash Store`n⚡️ Up to 5% off`n👉🏻 amzaff.in/l3swo0g`n🌟 Ki"
$ws.Range("D109").Value = "Scheduled Image + Text at 2025-07-28 12:00:00"
$ws.Range("E109").Value = "2025-07-28 12:00:05"
$ws.Range("F109").Value = "Kid's Carnival"
